# "At last, TSS proof finished"
# Fix the pseudocode on slide 5 ("TS Stack Pop"):
#   - "if(maxTS < n->TS){"    ->  "if(maxTS < n->ts){"
#   - "      youngest = ts;"  ->  "      youngest = n;"

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)          # "Content Placeholder 2" - the code listing
$tr  = $shp.TextFrame.TextRange

# --- "if(maxTS < n->TS){" -> "if(maxTS < n->ts){" -------------------------
$f1 = $tr.Find("maxTS < n->TS){")

# Re-type the ">" in place so it becomes its own run (matches the
# fine-grained run split produced by the original edit).
$gt = $tr.Characters($f1.Start + 10, 1)
$gt.Text = ">"

# Lower-case the trailing "TS" to "ts".
$ts1 = $tr.Characters($f1.Start + 11, 2)
$ts1.Text = "ts"

# --- "youngest = ts;" -> "youngest = n;" -----------------------------------
$f2 = $tr.Find("youngest = ts;")
$n2 = $tr.Characters($f2.Start + 11, 2)
$n2.Text = "n"
